$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Value = "'" + $val
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
}

# Row 18
$ws.Range("A18").Value = " Abu Dhabi"
$ws.Range("B18").Value = " October 25 2020"
$ws.Range("C18").Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Range("D18").Value = "Mumbai Indians"
$ws.Range("E18").Value = "Rajasthan Royals"
$ws.Range("F18").Value = "Quinton de Kock †"
Set-TextValue "G18" "6"
Set-TextValue "H18" "4"
Set-TextValue "I18" "0"
Set-TextValue "J18" "1"
Set-TextValue "K18" "150.00"

# Row 19
$ws.Range("A19").Value = " Abu Dhabi"
$ws.Range("B19").Value = " September 19 2020"
$ws.Range("C19").Value = "Super Kings won by 5 wickets (with 4 balls remaining)"
$ws.Range("D19").Value = "Mumbai Indians"
$ws.Range("E19").Value = "Chennai Super Kings"
$ws.Range("F19").Value = "Quinton de Kock †"
Set-TextValue "G19" "33"
Set-TextValue "H19" "20"
Set-TextValue "I19" "5"
Set-TextValue "J19" "0"
Set-TextValue "K19" "165.00"

# Row 20
$ws.Range("A20").Value = " Abu Dhabi"
$ws.Range("B20").Value = " October 28 2020"
$ws.Range("C20").Value = "Mumbai won by 5 wickets (with 5 balls remaining)"
$ws.Range("D20").Value = "Mumbai Indians"
$ws.Range("E20").Value = "Royal Challengers Bangalore"
$ws.Range("F20").Value = "Quinton de Kock †"
Set-TextValue "G20" "18"
Set-TextValue "H20" "19"
Set-TextValue "I20" "0"
Set-TextValue "J20" "1"
Set-TextValue "K20" "94.73"

# Row 21
$ws.Range("A21").Value = " Sharjah"
$ws.Range("B21").Value = " October 23 2020"
$ws.Range("C21").Value = "Mumbai won by 10 wickets (with 46 balls remaining)"
$ws.Range("D21").Value = "Mumbai Indians"
$ws.Range("E21").Value = "Chennai Super Kings"
$ws.Range("F21").Value = "Quinton de Kock †"
Set-TextValue "G21" "46"
Set-TextValue "H21" "37"
Set-TextValue "I21" "5"
Set-TextValue "J21" "2"
Set-TextValue "K21" "124.32"

# Row 22
$ws.Range("A22").Value = " Abu Dhabi"
$ws.Range("B22").Value = " September 23 2020"
$ws.Range("C22").Value = "Mumbai won by 49 runs"
$ws.Range("D22").Value = "Mumbai Indians"
$ws.Range("E22").Value = "Kolkata Knight Riders"
$ws.Range("F22").Value = "Quinton de Kock †"
Set-TextValue "G22" "1"
Set-TextValue "H22" "3"
Set-TextValue "I22" "0"
Set-TextValue "J22" "0"
Set-TextValue "K22" "33.33"

# Row 23
$ws.Range("A23").Value = " Abu Dhabi"
$ws.Range("B23").Value = " October 11 2020"
$ws.Range("C23").Value = "Mumbai won by 5 wickets (with 2 balls remaining)"
$ws.Range("D23").Value = "Mumbai Indians"
$ws.Range("E23").Value = "Delhi Capitals"
$ws.Range("F23").Value = "Quinton de Kock †"
Set-TextValue "G23" "53"
Set-TextValue "H23" "36"
Set-TextValue "I23" "4"
Set-TextValue "J23" "3"
Set-TextValue "K23" "147.22"

# Row 24
$ws.Range("A24").Value = " Dubai (DSC)"
$ws.Range("B24").Value = " November 05 2020"
$ws.Range("C24").Value = "Mumbai won by 57 runs"
$ws.Range("D24").Value = "Mumbai Indians"
$ws.Range("E24").Value = "Delhi Capitals"
$ws.Range("F24").Value = "Quinton de Kock †"
Set-TextValue "G24" "40"
Set-TextValue "H24" "25"
Set-TextValue "I24" "5"
Set-TextValue "J24" "1"
Set-TextValue "K24" "160.00"

# Row 25
$ws.Range("A25").Value = " Abu Dhabi"
$ws.Range("B25").Value = " October 01 2020"
$ws.Range("C25").Value = "Mumbai won by 48 runs"
$ws.Range("D25").Value = "Mumbai Indians"
$ws.Range("E25").Value = "Kings XI Punjab"
$ws.Range("F25").Value = "Quinton de Kock †"
Set-TextValue "G25" "0"
Set-TextValue "H25" "5"
Set-TextValue "I25" "0"
Set-TextValue "J25" "0"
Set-TextValue "K25" "0.00"

# Row 26
$ws.Range("A26").Value = " Sharjah"
$ws.Range("B26").Value = " October 04 2020"
$ws.Range("C26").Value = "Mumbai won by 34 runs"
$ws.Range("D26").Value = "Mumbai Indians"
$ws.Range("E26").Value = "Sunrisers Hyderabad"
$ws.Range("F26").Value = "Quinton de Kock †"
Set-TextValue "G26" "67"
Set-TextValue "H26" "39"
Set-TextValue "I26" "4"
Set-TextValue "J26" "4"
Set-TextValue "K26" "171.79"

# Row 27
$ws.Range("A27").Value = " Sharjah"
$ws.Range("B27").Value = " November 03 2020"
$ws.Range("C27").Value = "Sunrisers won by 10 wickets (with 17 balls remaining)"
$ws.Range("D27").Value = "Mumbai Indians"
$ws.Range("E27").Value = "Sunrisers Hyderabad"
$ws.Range("F27").Value = "Quinton de Kock †"
Set-TextValue "G27" "25"
Set-TextValue "H27" "13"
Set-TextValue "I27" "2"
Set-TextValue "J27" "2"
Set-TextValue "K27" "192.30"

# Row 28
$ws.Range("A28").Value = " Abu Dhabi"
$ws.Range("B28").Value = " October 16 2020"
$ws.Range("C28").Value = "Mumbai won by 8 wickets (with 19 balls remaining)"
$ws.Range("D28").Value = "Mumbai Indians"
$ws.Range("E28").Value = "Kolkata Knight Riders"
$ws.Range("F28").Value = "Quinton de Kock †"
Set-TextValue "G28" "78"
Set-TextValue "H28" "44"
Set-TextValue "I28" "9"
Set-TextValue "J28" "3"
Set-TextValue "K28" "177.27"

# Row 29
$ws.Range("A29").Value = " Dubai (DSC)"
$ws.Range("B29").Value = " October 31 2020"
$ws.Range("C29").Value = "Mumbai won by 9 wickets (with 34 balls remaining)"
$ws.Range("D29").Value = "Mumbai Indians"
$ws.Range("E29").Value = "Delhi Capitals"
$ws.Range("F29").Value = "Quinton de Kock †"
Set-TextValue "G29" "26"
Set-TextValue "H29" "28"
Set-TextValue "I29" "2"
Set-TextValue "J29" "0"
Set-TextValue "K29" "92.85"

# Row 30
$ws.Range("A30").Value = " Dubai (DSC)"
$ws.Range("B30").Value = " November 10 2020"
$ws.Range("C30").Value = "Mumbai won by 5 wickets (with 8 balls remaining)"
$ws.Range("D30").Value = "Mumbai Indians"
$ws.Range("E30").Value = "Delhi Capitals"
$ws.Range("F30").Value = "Quinton de Kock †"
Set-TextValue "G30" "20"
Set-TextValue "H30" "12"
Set-TextValue "I30" "3"
Set-TextValue "J30" "1"
Set-TextValue "K30" "166.66"

# Row 31
$ws.Range("A31").Value = " Dubai (DSC)"
$ws.Range("B31").Value = " September 28 2020"
$ws.Range("C31").Value = "Match tied (RCB won the one-over eliminator)"
$ws.Range("D31").Value = "Mumbai Indians"
$ws.Range("E31").Value = "Royal Challengers Bangalore"
$ws.Range("F31").Value = "Quinton de Kock †"
Set-TextValue "G31" "14"
Set-TextValue "H31" "15"
Set-TextValue "I31" "1"
Set-TextValue "J31" "0"
Set-TextValue "K31" "93.33"

# Row 32
$ws.Range("A32").Value = " Abu Dhabi"
$ws.Range("B32").Value = " October 06 2020"
$ws.Range("C32").Value = "Mumbai won by 57 runs"
$ws.Range("D32").Value = "Mumbai Indians"
$ws.Range("E32").Value = "Rajasthan Royals"
$ws.Range("F32").Value = "Quinton de Kock †"
Set-TextValue "G32" "23"
Set-TextValue "H32" "15"
Set-TextValue "I32" "3"
Set-TextValue "J32" "1"
Set-TextValue "K32" "153.33"

# Row 33
$ws.Range("A33").Value = " Dubai (DSC)"
$ws.Range("B33").Value = " October 18 2020"
$ws.Range("C33").Value = "Match tied (Kings XI won the one-over eliminator)"
$ws.Range("D33").Value = "Mumbai Indians"
$ws.Range("E33").Value = "Kings XI Punjab"
$ws.Range("F33").Value = "Quinton de Kock †"
Set-TextValue "G33" "53"
Set-TextValue "H33" "43"
Set-TextValue "I33" "3"
Set-TextValue "J33" "3"
Set-TextValue "K33" "123.25"

